# Refresh the cryptos list: update Price (D) / Volume(1h) (E) columns
# with the latest scrape, and fix rows 40-41 where
# PolygonEcosystemToken and WhiteBITCoin had been swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.372.94"
$ws.Range("E2").Value = "  +3.93%  "

$ws.Range("D3").Value = "3.128.03"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "205.12"
$ws.Range("E5").Value = "  +3.02%  "

$ws.Range("D6").Value = "620.45"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").Value = "0.263"
$ws.Range("E7").Value = "  +23.98%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +4.98%  "

$ws.Range("D10").Value = "3.125.18"
$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("D11").Value = "0.579"
$ws.Range("E11").Value = "  +29.82%  "

$ws.Range("D12").Value = "0.0000245"
$ws.Range("E12").Value = "  +25.30%  "

$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").Value = "3.698.98"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").Value = "31.13"
$ws.Range("E16").Value = "  +6.06%  "

$ws.Range("D17").Value = "79.131.02"
$ws.Range("E17").Value = "  +3.87%  "

$ws.Range("D18").Value = "3.121.34"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("D19").Value = "14.03"
$ws.Range("E19").Value = "  +3.56%  "

$ws.Range("D20").Value = "2.92"
$ws.Range("E20").Value = "  +12.69%  "

$ws.Range("D21").Value = "429.11"
$ws.Range("E21").Value = "  +12.12%  "

$ws.Range("D22").Value = "8.97"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").Value = "5.14"
$ws.Range("E23").Value = "  +14.28%  "

$ws.Range("D24").Value = "6.78"
$ws.Range("E24").Value = "  +5.22%  "

$ws.Range("D25").Value = "3.286.87"
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("E26").Value = "  +3.35%  "

$ws.Range("D27").Value = "4.61"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").Value = "10.67"
$ws.Range("E28").Value = "  +6.09%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  +10.01%  "

$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("D32").Value = "8.84"
$ws.Range("E32").Value = "  +6.18%  "

$ws.Range("D33").Value = "548.57"
$ws.Range("E33").Value = "  +9.18%  "

$ws.Range("E34").Value = "  +1.15%  "

$ws.Range("D35").Value = "0.149"
$ws.Range("E35").Value = "  +18.18%  "

$ws.Range("E36").Value = "  +1.79%  "

$ws.Range("D37").Value = "22.59"
$ws.Range("E37").Value = "  +8.23%  "

$ws.Range("E38").Value = "  +18.07%  "

$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "20.72"
$ws.Range("E40").Value = "  +3.21%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "0.397"
$ws.Range("E41").Value = "  +4.71%  "

$ws.Range("D42").Value = "162.18"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").Value = "5.48"
$ws.Range("E44").Value = "  +6.03%  "

$ws.Range("D45").Value = "'186.90"
$ws.Range("E45").Value = "  -4.34%  "

$ws.Range("E46").Value = "  +6.60%  "

$ws.Range("E47").Value = "  +8.03%  "

$ws.Range("E48").Value = "  -5.18%  "

$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("D50").Value = "42.49"
$ws.Range("E50").Value = "  +4.32%  "

$ws.Range("D51").Value = "4.17"
